$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.130.18"
$ws.Cells.Item(2, 5).Value = "  +5.48%  "
$ws.Cells.Item(3, 4).Value = "1.920.59"
$ws.Cells.Item(3, 5).Value = "  +2.21%  "
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.99%  "
$ws.Cells.Item(5, 4).Value = "326.99"
$ws.Cells.Item(5, 5).Value = "  +3.03%  "
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(7, 4).Value = "0.5167"
$ws.Cells.Item(7, 5).Value = "  +1.60%  "
$ws.Cells.Item(8, 4).Value = "0.4012"
$ws.Cells.Item(8, 5).Value = "  +2.85%  "
$ws.Cells.Item(9, 4).Value = "0.08455"
$ws.Cells.Item(9, 5).Value = "  +0.53%  "
$ws.Cells.Item(10, 2).Value = "OKB"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(10, 4).Value = "42.79"
$ws.Cells.Item(10, 5).Value = "  +2.32%  "
$ws.Cells.Item(11, 2).Value = "Polygon"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(11, 4).Value = "1.122"
$ws.Cells.Item(11, 5).Value = "  +1.54%  "
$ws.Cells.Item(12, 4).Value = "22.20"
$ws.Cells.Item(12, 5).Value = "  +8.67%  "
$ws.Cells.Item(13, 4).Value = "6.327"
$ws.Cells.Item(13, 5).Value = "  +1.56%  "
$ws.Cells.Item(14, 4).Value = "1.921.53"
$ws.Cells.Item(14, 5).Value = "  +2.40%  "
$ws.Cells.Item(15, 4).Value = "7.353"
$ws.Cells.Item(15, 5).Value = "  +1.43%  "
$ws.Cells.Item(16, 4).Value = "1.001"
$ws.Cells.Item(16, 5).Value = "  -1.34%  "
$ws.Cells.Item(17, 4).Value = "96.23"
$ws.Cells.Item(17, 5).Value = "  +5.35%  "
$ws.Cells.Item(18, 5).Value = "  +1.00%  "
$ws.Cells.Item(19, 4).Value = "0.06717"
$ws.Cells.Item(19, 5).Value = "  -0.29%  "
$ws.Cells.Item(20, 4).Value = "18.16"
$ws.Cells.Item(20, 5).Value = "  +2.40%  "
$ws.Cells.Item(21, 5).Value = "  -0.86%  "
$ws.Cells.Item(22, 5).Value = "  +1.92%  "
$ws.Cells.Item(23, 4).Value = "30.134.73"
$ws.Cells.Item(23, 5).Value = "  +5.30%  "
$ws.Cells.Item(24, 4).Value = "11.24"
$ws.Cells.Item(24, 5).Value = "  +1.19%  "
$ws.Cells.Item(25, 5).Value = "  -1.68%  "
$ws.Cells.Item(26, 4).Value = "2.144.63"
$ws.Cells.Item(26, 5).Value = "  +2.66%  "
$ws.Cells.Item(27, 4).Value = "160.60"
$ws.Cells.Item(27, 5).Value = "  -0.96%  "
$ws.Cells.Item(28, 4).Value = "21.03"
$ws.Cells.Item(28, 5).Value = "  +1.69%  "
$ws.Cells.Item(29, 4).Value = "2.460"
$ws.Cells.Item(29, 5).Value = "  +3.35%  "
$ws.Cells.Item(30, 4).Value = "128.88"
$ws.Cells.Item(30, 5).Value = "  +2.16%  "
$ws.Cells.Item(31, 4).Value = "1.074"
$ws.Cells.Item(31, 5).Value = "  +3.13%  "
$ws.Cells.Item(32, 4).Value = "0.1058"
$ws.Cells.Item(32, 5).Value = "  +1.16%  "
$ws.Cells.Item(33, 4).Value = "6.077"
$ws.Cells.Item(33, 5).Value = "  +4.94%  "
$ws.Cells.Item(34, 4).Value = "3.665"
$ws.Cells.Item(34, 5).Value = "  +1.09%  "
$ws.Cells.Item(35, 4).Value = "0.02519"
$ws.Cells.Item(35, 5).Value = "  +2.13%  "
$ws.Cells.Item(36, 4).Value = "0.06587"
$ws.Cells.Item(36, 5).Value = "  +0.47%  "
$ws.Cells.Item(37, 4).Value = "0.2216"
$ws.Cells.Item(37, 5).Value = "  +2.48%  "
$ws.Cells.Item(38, 5).Value = "  +3.59%  "
$ws.Cells.Item(39, 4).Value = "9.005"
$ws.Cells.Item(39, 5).Value = "  +1.37%  "
$ws.Cells.Item(40, 4).Value = "5.204"
$ws.Cells.Item(40, 5).Value = "  +2.22%  "
$ws.Cells.Item(41, 4).Value = "0.6548"
$ws.Cells.Item(41, 5).Value = "  +1.80%  "
$ws.Cells.Item(42, 4).Value = "1.246"
$ws.Cells.Item(42, 5).Value = "  -0.80%  "
$ws.Cells.Item(43, 4).Value = "11.43"
$ws.Cells.Item(43, 5).Value = "  +2.66%  "
$ws.Cells.Item(44, 4).Value = "0.6145"
$ws.Cells.Item(44, 5).Value = "  +1.71%  "
$ws.Cells.Item(45, 4).Value = "13.18"
$ws.Cells.Item(45, 5).Value = "  +1.23%  "
$ws.Cells.Item(46, 4).Value = "3.755"
$ws.Cells.Item(46, 5).Value = "  +1.44%  "
$ws.Cells.Item(47, 4).Value = "2.057"
$ws.Cells.Item(47, 5).Value = "  +2.28%  "
$ws.Cells.Item(48, 4).Value = "1.244"
$ws.Cells.Item(48, 5).Value = "  +2.01%  "
$ws.Cells.Item(49, 4).Value = "125.29"
$ws.Cells.Item(49, 5).Value = "  +2.58%  "
$ws.Cells.Item(50, 5).Value = "  +2.54%  "
$ws.Cells.Item(51, 4).Value = "79.24"
$ws.Cells.Item(51, 5).Value = "  +3.25%  "
